# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates/creations/deletions to match the target diff
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 925
$ws.Range("I12").Value = 699.6
$ws.Range("J12").Value = 1300.6666
$ws.Range("K12").Value = 699.6
$ws.Range("L12").Value = 1300.6666
$ws.Range("M12").Value = -529.6
$ws.Range("N12").Value = -1640.6666
$ws.Range("H29").Value = 1347.6154
$ws.Range("I29").Value = 524.8
$ws.Range("K29").Value = 1574.4
$ws.Range("M29").Value = -1293.4
$ws.Range("H51").Value = 8636.362999999999
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 8636.362999999999
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 8636.362999999999
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -9604.362999999999
$ws.Range("H64").Value = 7882.846
$ws.Range("I64").Value = 7135.2
$ws.Range("J64").Value = 8350.125
$ws.Range("K64").Value = 7135.2
$ws.Range("L64").Value = 8350.125
$ws.Range("M64").Value = -6887.2
$ws.Range("N64").Value = -8846.125
$ws.Range("H67").Value = 7882.846
$ws.Range("I67").Value = 7135.2
$ws.Range("J67").Value = 8350.125
$ws.Range("K67").Value = 7135.2
$ws.Range("L67").Value = 8350.125
$ws.Range("M67").Value = -6277.2
$ws.Range("N67").Value = -10066.125
$ws.Range("H74").Value = 11454.818
$ws.Range("I74").Value = 9999.75
$ws.Range("K74").Value = 9999.75
$ws.Range("M74").Value = -9063.75
$ws.Range("H77").Value = 11454.818
$ws.Range("I77").Value = 9999.75
$ws.Range("K77").Value = 49998.75
$ws.Range("M77").Value = -45318.75
$ws.Range("H112").Value = 1441.4546
$ws.Range("J112").Value = 1524.1111
$ws.Range("L112").Value = 4572.3333
$ws.Range("N112").Value = -6788.3333
$ws.Range("H113").Value = 4908.636
$ws.Range("I113").Value = 1599
$ws.Range("K113").Value = 1599
$ws.Range("M113").Value = 1655
$ws.Range("H125").Value = 2825.111
$ws.Range("I125").Value = 2463.3333
$ws.Range("K125").Value = 22169.9997
$ws.Range("M125").Value = -19709.9997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5468.2
$ws.Range("I2").Value = 365.52942
$ws.Range("K2").Value = 365.52942
$ws.Range("M2").Value = -252.52942
$ws.Range("H32").Value = 6547
$ws.Range("I32").Value = 5101
$ws.Range("K32").Value = 5101
$ws.Range("M32").Value = -4814
$ws.Range("H110").Value = 2655.1072
$ws.Range("I110").Value = 1873.2
$ws.Range("J110").Value = 9171
$ws.Range("K110").Value = 1873.2
$ws.Range("L110").Value = 9171
$ws.Range("M110").Value = 171.8
$ws.Range("N110").Value = -13261
$ws.Range("H112").Value = 43786.285
$ws.Range("J112").Value = 43786.285
$ws.Range("L112").Value = 43786.285
$ws.Range("N112").Value = -46740.285
$ws.Range("H116").Value = 5468.2
$ws.Range("I116").Value = 365.52942
$ws.Range("K116").Value = 365.52942
$ws.Range("M116").Value = 1928.47058
$ws.Range("H122").Value = 2929.9666
$ws.Range("I122").Value = 2297.1667
$ws.Range("K122").Value = 6891.500100000001
$ws.Range("M122").Value = -4441.500100000001
$ws.Range("H132").Value = 2278.4546
$ws.Range("I132").Value = 1650.8064
$ws.Range("K132").Value = 4952.4192
$ws.Range("M132").Value = -2422.4192
$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5468.2
$ws.Range("I3").Value = 365.52942
$ws.Range("K3").Value = 365.52942
$ws.Range("M3").Value = -251.52942

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3985
$ws.Range("I105").Value = 1576
$ws.Range("J105").Value = 6996.25
$ws.Range("K105").Value = 1576
$ws.Range("L105").Value = 6996.25
$ws.Range("M105").Value = 171
$ws.Range("N105").Value = -10490.25
$ws.Range("H141").Value = 274042.78
$ws.Range("J141").Value = 274042.78
$ws.Range("L141").Value = 274042.78
$ws.Range("N141").Value = -284402.78

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 179.45454
$ws.Range("I7").Value = 185
$ws.Range("J7").Value = 172.8
$ws.Range("K7").Value = 555
$ws.Range("L7").Value = 518.4000000000001
$ws.Range("M7").Value = -443
$ws.Range("N7").Value = -742.4000000000001
$ws.Range("H107").Value = 66667836
$ws.Range("I107").Value = 791.6667
$ws.Range("J107").Value = 83334600
$ws.Range("K107").Value = 2375.0001
$ws.Range("L107").Value = 250003800
$ws.Range("M107").Value = -455.0001000000002
$ws.Range("N107").Value = -250007640
$ws.Range("H131").Value = 7987934.5
$ws.Range("J131").Value = 5557472
$ws.Range("L131").Value = 16672416
$ws.Range("N131").Value = -16682496

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1778.4445
$ws.Range("I97").Value = 1466.7333
$ws.Range("J97").Value = 3337
$ws.Range("K97").Value = 1466.7333
$ws.Range("L97").Value = 3337
$ws.Range("M97").Value = -970.7333000000001
$ws.Range("N97").Value = -4329
$ws.Range("H102").Value = 2206.3684
$ws.Range("I102").Value = 1492.4884
$ws.Range("J102").Value = 4399
$ws.Range("K102").Value = 1492.4884
$ws.Range("L102").Value = 4399
$ws.Range("M102").Value = 129.5116
$ws.Range("N102").Value = -7643
$ws.Range("H126").Value = 4659.0557
$ws.Range("I126").Value = 3057.625
$ws.Range("J126").Value = 5940.2
$ws.Range("K126").Value = 9172.875
$ws.Range("L126").Value = 17820.6
$ws.Range("M126").Value = -6702.875
$ws.Range("N126").Value = -22760.6
$ws.Range("H132").Value = 4980.3335
$ws.Range("I132").Value = 4703.3335
$ws.Range("K132").Value = 14110.0005
$ws.Range("M132").Value = -11580.0005
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1614356.1
$ws.Range("I55").Value = 2631961.2
$ws.Range("J55").Value = 3148
$ws.Range("K55").Value = 2631961.2
$ws.Range("L55").Value = 3148
$ws.Range("M55").Value = -2631788.2
$ws.Range("N55").Value = -3494
$ws.Range("H82").Value = 25268.965
$ws.Range("I82").Value = 18103.268
$ws.Range("K82").Value = 18103.268
$ws.Range("M82").Value = -17742.268
$ws.Range("H85").Value = 25268.965
$ws.Range("I85").Value = 18103.268
$ws.Range("K85").Value = 18103.268
$ws.Range("M85").Value = -16855.268
$ws.Range("H93").Value = 2023.75
$ws.Range("I93").Value = 2023.75
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2023.75
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -775.75
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 75000000
$ws.Range("J5").Value = 75000000
$ws.Range("L5").Value = 75000000
$ws.Range("N5").Value = -75000224
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("H86").Value = 44497.5
$ws.Range("I86").Value = 39995
$ws.Range("K86").Value = 39995
$ws.Range("M86").Value = -38872
$ws.Range("H89").Value = 44497.5
$ws.Range("I89").Value = 39995
$ws.Range("K89").Value = 199975
$ws.Range("M89").Value = -194359
$ws.Range("H113").Value = 465.97437
$ws.Range("I113").Value = 356.64
$ws.Range("K113").Value = 1069.92
$ws.Range("M113").Value = 1100.08
$ws.Range("H122").Value = 1730.2333
$ws.Range("I122").Value = 768.5714
$ws.Range("K122").Value = 2305.7142
$ws.Range("M122").Value = 144.2857999999997

